$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.226.46"
$ws.Range("E2").Value = "  -3.52%  "
# Row 3
$ws.Range("D3").Value = "2.246.66"
$ws.Range("E3").Value = "  -4.25%  "
# Row 4
$ws.Range("E4").Value = "  +0.01%  "
# Row 5
$ws.Range("D5").Value = "'244.65"
$ws.Range("E5").Value = "  +2.29%  "
# Row 6
$ws.Range("D6").Value = "'0.632"
$ws.Range("E6").Value = "  -4.96%  "
# Row 7
$ws.Range("D7").Value = "'69.41"
$ws.Range("E7").Value = "  -4.11%  "
# Row 8
$ws.Range("E8").Value = "  +0.12%  "
# Row 9
$ws.Range("D9").Value = "'0.555"
$ws.Range("E9").Value = "  -6.17%  "
# Row 10
$ws.Range("D10").Value = "'0.0990"
$ws.Range("E10").Value = "  -0.98%  "
# Row 11
$ws.Range("D11").Value = "'58.99"
$ws.Range("E11").Value = "  +1.10%  "
# Row 12
$ws.Range("D12").Value = "'36.39"
# Row 13
$ws.Range("E13").Value = "  -1.61%  "
# Row 14
$ws.Range("D14").Value = "'6.82"
$ws.Range("E14").Value = "  -5.57%  "
# Row 15
$ws.Range("D15").Value = "2.580.23"
$ws.Range("E15").Value = "  -4.24%  "
# Row 16
$ws.Range("D16").Value = "'15.03"
$ws.Range("E16").Value = "  -6.39%  "
# Row 17
$ws.Range("D17").Value = "'0.868"
$ws.Range("E17").Value = "  -3.48%  "
# Row 18
$ws.Range("D18").Value = "2.254.19"
$ws.Range("E18").Value = "  -3.71%  "
# Row 19
$ws.Range("D19").Value = "42.222.86"
$ws.Range("E19").Value = "  -3.39%  "
# Row 20
$ws.Range("D20").Value = "0.0₃0972"
$ws.Range("E20").Value = "  -3.86%  "
# Row 21
$ws.Range("D21").Value = "'6.27"
$ws.Range("E21").Value = "  -5.35%  "
# Row 22
$ws.Range("D22").Value = "'73.42"
$ws.Range("E22").Value = "  -5.48%  "
# Row 23
$ws.Range("D23").Value = "'237.09"
$ws.Range("E23").Value = "  -5.54%  "
# Row 24
$ws.Range("D24").Value = "'2.04"
$ws.Range("E24").Value = "  +10.27%  "
# Row 25
$ws.Range("E25").Value = "  -0.18%  "
# Row 26
$ws.Range("E26").Value = "  -1.26%  "
# Row 27
$ws.Range("D27").Value = "'2.48"
$ws.Range("E27").Value = "  -0.01%  "
# Row 28
$ws.Range("D28").Value = "'10.05"
$ws.Range("E28").Value = "  -2.83%  "
# Row 29
$ws.Range("D29").Value = "'2.24"
$ws.Range("E29").Value = "  -1.33%  "
# Row 30
$ws.Range("D30").Value = "'172.99"
$ws.Range("E30").Value = "  -2.18%  "
# Row 31
$ws.Range("D31").Value = "'20.58"
$ws.Range("E31").Value = "  -7.26%  "
# Row 32
$ws.Range("D32").Value = "'0.123"
$ws.Range("E32").Value = "  -1.80%  "
# Row 34
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "'5.37"
$ws.Range("E34").Value = "  +1.11%  "
# Row 35
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.0720"
$ws.Range("E35").Value = "  -4.21%  "
# Row 36
$ws.Range("D36").Value = "'4.71"
$ws.Range("E36").Value = "  -6.91%  "
# Row 37
$ws.Range("D37").Value = "'3.87"
$ws.Range("E37").Value = "  +3.52%  "
# Row 38
$ws.Range("D38").Value = "'22.80"
$ws.Range("E38").Value = "  +20.90%  "
# Row 39
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.0283"
$ws.Range("E39").Value = "  +4.48%  "
# Row 40
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "'2.32"
$ws.Range("E40").Value = "  -1.35%  "
# Row 41
$ws.Range("D41").Value = "'5.93"
$ws.Range("E41").Value = "  -7.06%  "
# Row 42
$ws.Range("D42").Value = "'66.51"
$ws.Range("E42").Value = "  +1.04%  "
# Row 43
$ws.Range("D43").Value = "'9.41"
$ws.Range("E43").Value = "  +2.32%  "
# Row 44
$ws.Range("D44").Value = "'4.97"
$ws.Range("E44").Value = "  -12.26%  "
# Row 45
$ws.Range("E45").Value = "  -3.58%  "
# Row 46
$ws.Range("D46").Value = "'0.191"
$ws.Range("E46").Value = "  -2.14%  "
# Row 47
$ws.Range("D47").Value = "'4.61"
$ws.Range("E47").Value = "  +11.14%  "
# Row 48
$ws.Range("E48").Value = "  +0.04%  "
# Row 49
$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D49").Value = "'1.20"
$ws.Range("E49").Value = "  -1.89%  "
# Row 50
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "'2.36"
$ws.Range("E50").Value = "  -2.06%  "
# Row 51
$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D51").Value = "'1.11"
$ws.Range("E51").Value = "  -2.95%  "
